$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the formatting of column N (2020) into the new column O (2021)
# without minting new cell styles, then fill in the 2021 figures.
$ws.Range("N3:N5").Copy()
$ws.Range("O3:O5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("O3").Value = 2021
$ws.Range("O4").Value = 14
$ws.Range("O5").Value = 1252.8

$ws.Range("O9").Select()
